$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("N2").Value = 14035.95
$ws.Range("P2").Value = 13987.79
$ws.Range("Q2").Value = 9198.35
$ws.Range("R2").Value = 4251.8
$ws.Range("AG2").Value = 189407.06

# Row 3
$ws.Range("P3").Value = 13874
$ws.Range("Q3").Value = 11612
$ws.Range("R3").Value = 3198
$ws.Range("AG3").Value = 117315.7

# Row 4
$ws.Range("P4").Value = 3986
$ws.Range("Q4").Value = 4253.01
$ws.Range("R4").Value = 3270
$ws.Range("S4").Value = 3749.01
$ws.Range("AG4").Value = 60627.6

# Row 5
$ws.Range("O5").Value = 1726
$ws.Range("P5").Value = 2093
$ws.Range("Q5").Value = 5340.4
$ws.Range("R5").Value = 3814
$ws.Range("S5").Value = 4336
$ws.Range("AG5").Value = 60516.19

# Row 6
$ws.Range("N6").Value = 26531.86
$ws.Range("O6").Value = 36574.18
$ws.Range("P6").Value = 33940.79
$ws.Range("Q6").Value = 30403.76
$ws.Range("R6").Value = 14533.8
$ws.Range("S6").Value = 8085.01
$ws.Range("AG6").Value = 427866.55

$wb.Save()
